$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_3_01")

# 1. Update the title (row 1) from "October 2016" to "November 2016"
$ws.Range("A1").Value = "Table 3.1. Stocks of Coal, Petroleum Liquids, and Petroleum Coke: Electric Power Sector, 2006 - November 2016"

# 2. Insert a new data row (row 52, "November") above the existing "Notes" row,
#    which pushes the Notes row (and its A52:J52 merge) down to row 53.
$ws.Rows.Item(52).Insert()

# Copy the formatting of the previous data row (51, "October") into the
# newly inserted row 52 so the new row matches the style of the other
# monthly data rows (style ids 9/10) instead of Excel's generic insert default.
$ws.Range("A51:J51").Copy()
$ws.Range("A52:J52").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new "November" row with its values
$ws.Range("A52").Value = "November"

$ws.Range("B52").Value = 172139
$ws.Range("C52").Value = 30847
$ws.Range("D52").Value = 833

$ws.Range("E52").Value = 139080
$ws.Range("F52").Value = 20372
$ws.Range("G52").Value = 606

$ws.Range("H52").Value = 33059
$ws.Range("I52").Value = 10475
$ws.Range("J52").Value = 227
